# Trade #65 closed at 2026-02-17 08:48:55 - unknown UNKNOWN +0.000%
#
# Appends the newly-closed trade (#65) to the "All Trades" and
# "MarketMaking" logs, and rolls the aggregate stats on the "Summary"
# and "Strategy Status" sheets forward to include it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value2 = 1199.62   # Current Capital
$summary.Range("B4").Value2 = -0.38     # Total P&L $
$summary.Range("B5").Value2 = -0.12     # Total P&L %
$summary.Range("B6").Value2 = 65        # Total Trades
$summary.Range("B7").Value2 = 26        # Winning Trades
$summary.Range("B9").Value2 = 40        # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value2 = 99.62      # Capital
$status.Range("D4").Value2 = 65         # Trades
$status.Range("E4").Value2 = -0.38      # P&L $
$status.Range("F4").Value2 = -0.38      # P&L %
$status.Range("G4").Value2 = 40         # Win Rate %

# ---------------------------------------------------------------------
# New trade row (#65) appended to both trade-log sheets
# ---------------------------------------------------------------------
function Add-Trade65Row($ws) {
    $row = 66

    $ws.Cells.Item($row, 1).Value2 = 65

    # Date column looks like "2026-02-17" which Excel's auto-detection
    # would otherwise coerce into a date serial number - force it to be
    # stored as text, matching the rest of the column.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value2 = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value2 = "08:48:49"
    $ws.Cells.Item($row, 4).Value2 = "MarketMaking"
    $ws.Cells.Item($row, 5).Value2 = "UP"
    $ws.Cells.Item($row, 6).Value2 = 0.88
    $ws.Cells.Item($row, 7).Value2 = 0.91
    $ws.Cells.Item($row, 8).Value2 = "CLOSED"
    $ws.Cells.Item($row, 9).Value2 = 3.4091
    $ws.Cells.Item($row, 10).Value2 = 0.03
    $ws.Cells.Item($row, 11).Value2 = 99.62
    $ws.Cells.Item($row, 12).Value2 = 0
    $ws.Cells.Item($row, 13).Value2 = 0
    $ws.Cells.Item($row, 14).Value2 = 0.6
    $ws.Cells.Item($row, 15).Value2 = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value2 = "early_exit"
    $ws.Cells.Item($row, 17).Value2 = 0.13
}

Add-Trade65Row $wb.Worksheets.Item("All Trades")
Add-Trade65Row $wb.Worksheets.Item("MarketMaking")
